# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Mon Sep  4 03:45:07 UTC 2023 with GitHub Actions"
#   - Price (D) / Volume(1h) (E) columns refreshed with new scrape values
#   - Rows 45 & 46 (BabyDogeCoin / RocketPoolETH) swapped position
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 45 / 46 swap (coin name + link) ---
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"

# --- Price (D) / Volume 1h (E) refresh ---
# A few Price values are plain decimals ("1.003", "216.38", ...) that Excel
# would otherwise auto-convert to numbers. Briefly force Text format so the
# literal string (incl. trailing zeros) is preserved, then restore General
# formatting to match the sheets original (unstyled) text cells.
$ws.Range("D2").Value = "26.026.60"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.643.23"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.38"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5160"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2598"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06395"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.92"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.313"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "1.653.38"
$ws.Range("E13").Value = "  -3.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5500"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.85"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "0.0₅7786"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").Value = "26.069.59"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "199.68"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.487"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +1.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.03"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.136"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +1.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("E24").Value = "  +2.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.16"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1221"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +6.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.912"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.74"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04892"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -3.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.313"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.251"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.546"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.385"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9211"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +2.97%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5604"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("D38").Value = "1.116.29"
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01576"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.547"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.575"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8122"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.80"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "1.782.29"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "0.0₈119"
$ws.Range("E46").Value = "  -2.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4541"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.41"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05258"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +3.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09617"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.54%  "
